$wb = $excel.ActiveWorkbook

# --- Rename the "gamelog" sheet to "Gamelog" ---
$wsGamelog = $wb.Worksheets.Item("gamelog")
$wsGamelog.Name = "Gamelog"

# --- Car sheet: update the remembered selection (no data change) ---
$wsCar = $wb.Worksheets.Item("Car")
$wsCar.Range("E24").Select()

# --- Gamelog sheet: restructure header row from 7 columns to 5 columns ---
# old: id, house_id, spouse_id, car_id, user_id, created_date, updated_date
# new: id, house, spouse, car, user_id
$wsGamelog.Range("B1").Value = "house"
$wsGamelog.Range("C1").Value = "spouse"
$wsGamelog.Range("D1").Value = "car"
$wsGamelog.Range("E1").Value = "user_id"
$wsGamelog.Range("F1:G1").Clear()

# Make Gamelog the active/selected sheet with the remembered selection
$wsGamelog.Range("F26").Select()
